$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet from the generic "1" to the descriptive "Zestaponi"
$ws.Name = "Zestaponi"

# The confidentiality/unavailability marker used throughout the sheet
$mark = "..."

# Row 6 ("Urban") - additional years become confidential/unavailable
$urbanCells = @("B6","D6","G6","I6","J6","K6","L6","M6","O6")
foreach ($addr in $urbanCells) {
    $ws.Range($addr).Value = $mark
}

# Row 7 ("Rural") - additional years become confidential/unavailable
$ruralCells = @("B7","D7","G7","J7","K7","L7","M7","N7","O7")
foreach ($addr in $ruralCells) {
    $ws.Range($addr).Value = $mark
}

# Row 5 ("Total") already has J5 marked; re-write it so the old ellipsis
# shared-string entry has no remaining references and gets dropped.
$ws.Range("J5").Value = $mark

# Delete the blank row 8 so the footnote row (old row 9) shifts up to row 8
$ws.Rows(8).Delete()
